# Updated symbol list (coin prices / 1h volume%, and a few re-ranked rows
# with their Coin name + Link) to match the latest coinranking.com pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "258.40"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.64%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.861"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-9.10%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05959"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.32%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.685"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.35%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8753"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.26%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9605"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.97%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1417"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.11%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07223"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.25%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03136"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.19%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09246"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.11%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001547"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0006068"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.06%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006025"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.22%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.483"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.47%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.223"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.68%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.219"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.32%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3144"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.70%"
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03574"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.99%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1308"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.19%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.530"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.08%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04246"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.88%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1380"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.10%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001220"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.31%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.35%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.05%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001493"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "2.67%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03833"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.50%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005888"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.78%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1105"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.30%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.60%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.90%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005489"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.35%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.06%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1091"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "9.17%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002151"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.71%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.06%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.06%"
